# Hortaliza, Feria Lagunitas de Puerto Montt - Cebolla
# Commit: "Fruta / hortaliza, semanal" -- weekly append of two new price
# records (dated 44946) inserted right before the existing row 703, pushing
# every subsequent record down by two rows (old 703..781 -> new 705..783).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the current row 703; each Insert() shifts the
# row that is already at 703 (and everything below it) down by one, so two
# calls give us the required two-row shift while copying the formatting of
# the row above (keeps the date-style on column D, etc.).
$ws.Rows.Item(703).Insert()
$ws.Rows.Item(703).Insert()

# --- New row 703 --------------------------------------------------------
$ws.Range("A703").Value = 4
$ws.Range("B703").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C703").Value = "Los Lagos"
$ws.Range("D703").Value = 44946
$ws.Range("E703").Value = 10
$ws.Range("F703").Value = 100112004
$ws.Range("G703").Value = "Cebolla"
$ws.Range("H703").Value = "Morada(o)"
$ws.Range("I703").Value = "1a (cosecha)"
$ws.Range("J703").Value = 250
$ws.Range("K703").Value = 15000
$ws.Range("L703").Value = 15000
$ws.Range("M703").Value = 15000
$ws.Range("N703").Value = "`$/malla 18 kilos"
$ws.Range("O703").Value = "Región de O'Higgins"
$ws.Range("P703").Value = 833
$ws.Range("Q703").Value = 18
$ws.Range("R703").Value = "Hortaliza"

# --- New row 704 --------------------------------------------------------
$ws.Range("A704").Value = 4
$ws.Range("B704").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C704").Value = "Los Lagos"
$ws.Range("D704").Value = 44946
$ws.Range("E704").Value = 10
$ws.Range("F704").Value = 100112004
$ws.Range("G704").Value = "Cebolla"
$ws.Range("H704").Value = "Sin especificar"
$ws.Range("I704").Value = "1a (cosecha)"
$ws.Range("J704").Value = 1200
$ws.Range("K704").Value = 15000
$ws.Range("L704").Value = 15000
$ws.Range("M704").Value = 15000
$ws.Range("N704").Value = "`$/malla 18 kilos"
$ws.Range("O704").Value = "Región de O'Higgins"
$ws.Range("P704").Value = 833
$ws.Range("Q704").Value = 18
$ws.Range("R704").Value = "Hortaliza"
